# Corrects the "indices" dataset: removes the "Fixed" category rows
# (Replacements, O&M, Land) which were erroneous, shifting the Input /
# Output / Metric rows up by three. Also re-points the active sheet/tab
# to "indices" with a selection, matching the authored workbook state.

$wb = $excel.ActiveWorkbook

$indices = $wb.Worksheets.Item("indices")

# Remove the three "Fixed" rows (Replacements / O&M / Land) -- rows 6-8
# on the "indices" sheet. Everything below shifts up automatically.
$indices.Rows("6:8").Delete()

# The authored workbook leaves "indices" as the active/selected sheet
# (rather than "designs"), with cell B19 selected.
$indices.Activate()
$indices.Range("B19").Select()
